$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

# Numeric cells - plain values
$ws.Cells.Item($row, 1).Value = 1581552000
$ws.Cells.Item($row, 5).Value = 0.25
$ws.Cells.Item($row, 6).Value = 0.275
$ws.Cells.Item($row, 7).Value = 0.25
$ws.Cells.Item($row, 8).Value = 0.265
$ws.Cells.Item($row, 9).Value = 27529300

# Text cells (date string, id, name) - must stay text, not get
# reinterpreted as a date serial / number, and must not keep a
# lingering explicit number-format style once written.
$textRange = $ws.Range("B17:D17")
$textRange.NumberFormat = "@"
$ws.Cells.Item($row, 2).Value = "2020-02-13"
$ws.Cells.Item($row, 3).Value = "0217"
$ws.Cells.Item($row, 4).Value = "PWRWELL"
$textRange.ClearFormats()

Write-Host "Row 17 added"
